$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.929.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.229.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.28%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +18.81%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.60%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.561.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.226.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.786.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.31%  "
$ws.Range("E25").Value = "  +10.21%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0729"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +24.42%  "
$ws.Range("E38").Value = "  +10.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.99%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +22.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.208"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.12%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +8.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.11%  "
